$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new (blank) rows at position 22. Rows 19-21 stay put; the old
# rows 22-24 get pushed down to 25-27; rows 22-24 become blank and need to
# be repopulated with the (renumbered) data that used to live in 19-21.
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(22).Insert()

# The newly-inserted rows 22-24 lost the per-column cell formatting that the
# rest of the data rows have (bordered/bold index column, datetime / date
# number formats). Copy that formatting down from rows 19-21 before filling
# in values, so the look matches the rest of the table.
$ws.Range("A19:M21").Copy()
$ws.Range("A22:M24").PasteSpecial(-4122)

# --- Row 19: brand-new measurement data (index 17) ---
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 45392.6875
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 6
$ws.Range("H19").Value = 3
$ws.Range("I19").Value = 52
$ws.Range("J19").Value = 131
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 45392
$ws.Range("M19").Value = 16

# --- Row 20: brand-new measurement data (index 18) ---
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 45392.69444444445
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 7
$ws.Range("I20").Value = 54
$ws.Range("J20").Value = 165
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 45392
$ws.Range("M20").Value = 16

# --- Row 21: brand-new measurement data (index 19) ---
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 45392.70138888889
$ws.Range("C21").Value = 4
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 4
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = 60
$ws.Range("J21").Value = 181
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 45392
$ws.Range("M21").Value = 16

# --- Row 22: former row-19 data, index renumbered 17 -> 20 ---
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = 45392.83333333334
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 2
$ws.Range("I22").Value = 34
$ws.Range("J22").Value = 79
$ws.Range("K22").Value = 11
$ws.Range("L22").Value = 45392
$ws.Range("M22").Value = 20

# --- Row 23: former row-20 data, index renumbered 18 -> 21 ---
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = 45392.84027777778
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 2
$ws.Range("I23").Value = 30
$ws.Range("J23").Value = 87
$ws.Range("K23").Value = 7
$ws.Range("L23").Value = 45392
$ws.Range("M23").Value = 20

# --- Row 24: former row-21 data, index renumbered 19 -> 22 ---
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = 45392.84722222222
$ws.Range("C24").Value = 2
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 1
$ws.Range("I24").Value = 48
$ws.Range("J24").Value = 81
$ws.Range("K24").Value = 5
$ws.Range("L24").Value = 45392
$ws.Range("M24").Value = 20

# --- Rows 25-27 already hold the old row 22-24 data after the shift; only
# the running-index column (A) needs to be renumbered by +3 ---
$ws.Range("A25").Value = 23
$ws.Range("A26").Value = 24
$ws.Range("A27").Value = 25
